$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 2290220.5
$ws.Range("I11").Value = 2290220.5
$ws.Range("K11").Value = 2290220.5
$ws.Range("M11").Value = -2290080.5
$ws.Range("H43").Value = 2144.1052
$ws.Range("I43").Value = 2750
$ws.Range("J43").Value = 831.3333
$ws.Range("K43").Value = 2750
$ws.Range("L43").Value = 831.3333
$ws.Range("M43").Value = -2681
$ws.Range("N43").Value = -969.3333
$ws.Range("H129").Value = 1409.1837
$ws.Range("J129").Value = 1320.762
$ws.Range("L129").Value = 3962.286
$ws.Range("N129").Value = -13962.286
$ws.Range("H137").Value = 3625.1724
$ws.Range("I137").Value = 1112.5
$ws.Range("J137").Value = 4027.2
$ws.Range("K137").Value = 3337.5
$ws.Range("L137").Value = 12081.6
$ws.Range("M137").Value = -787.5
$ws.Range("N137").Value = -17181.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4722.4517
$ws.Range("I2").Value = 5369.115
$ws.Range("K2").Value = 5369.115
$ws.Range("M2").Value = -5256.115
$ws.Range("H32").Value = 2707.54
$ws.Range("I32").Value = 2368.598
$ws.Range("K32").Value = 2368.598
$ws.Range("M32").Value = -2081.598
$ws.Range("H116").Value = 4722.4517
$ws.Range("I116").Value = 5369.115
$ws.Range("K116").Value = 5369.115
$ws.Range("M116").Value = -3075.115
$ws.Range("H132").Value = 10640192
$ws.Range("I132").Value = 14707153
$ws.Range("J132").Value = 3526.4614
$ws.Range("K132").Value = 44121459
$ws.Range("L132").Value = 10579.3842
$ws.Range("M132").Value = -44118929
$ws.Range("N132").Value = -15639.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4722.4517
$ws.Range("I3").Value = 5369.115
$ws.Range("K3").Value = 5369.115
$ws.Range("M3").Value = -5255.115
$ws.Range("H38").Value = 33625
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H105").Value = 1823.5
$ws.Range("I105").Value = 1867.8422
$ws.Range("J105").Value = 1703.1428
$ws.Range("K105").Value = 1867.8422
$ws.Range("L105").Value = 1703.1428
$ws.Range("M105").Value = -120.8422
$ws.Range("N105").Value = -5197.1428
$ws.Range("H107").Value = 2309.0588
$ws.Range("I107").Value = 2160.0715
$ws.Range("J107").Value = 3004.3333
$ws.Range("K107").Value = 2160.0715
$ws.Range("L107").Value = 3004.3333
$ws.Range("M107").Value = -240.0715
$ws.Range("N107").Value = -6844.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 228.15384
$ws.Range("I7").Value = 136.6875
$ws.Range("J7").Value = 374.5
$ws.Range("K7").Value = 136.6875
$ws.Range("L7").Value = 374.5
$ws.Range("M7").Value = -23.6875
$ws.Range("N7").Value = -600.5
$ws.Range("H31").Value = 4580.75
$ws.Range("I31").Value = 1565.9259
$ws.Range("J31").Value = 6566.122
$ws.Range("K31").Value = 1565.9259
$ws.Range("L31").Value = 6566.122
$ws.Range("M31").Value = -1270.9259
$ws.Range("N31").Value = -7156.122
$ws.Range("H34").Value = 4580.75
$ws.Range("I34").Value = 1565.9259
$ws.Range("J34").Value = 6566.122
$ws.Range("K34").Value = 1565.9259
$ws.Range("L34").Value = 6566.122
$ws.Range("M34").Value = -1363.9259
$ws.Range("N34").Value = -6970.122
$ws.Range("H107").Value = 1008.7143
$ws.Range("I107").Value = 952.2
$ws.Range("J107").Value = 1150
$ws.Range("K107").Value = 952.2
$ws.Range("L107").Value = 1150
$ws.Range("M107").Value = 967.8
$ws.Range("N107").Value = -4990
$ws.Range("H132").Value = 55005.777
$ws.Range("I132").Value = 1983.5625
$ws.Range("J132").Value = 132129
$ws.Range("K132").Value = 5950.6875
$ws.Range("L132").Value = 396387
$ws.Range("M132").Value = -3420.6875
$ws.Range("N132").Value = -401447

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3623.2122
$ws.Range("I113").Value = 4726.5835
$ws.Range("J113").Value = 680.8889
$ws.Range("K113").Value = 14179.7505
$ws.Range("L113").Value = 2042.6667
$ws.Range("M113").Value = -12009.7505
$ws.Range("N113").Value = -6382.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2995.6924
$ws.Range("I102").Value = 2941.6
$ws.Range("J102").Value = 3176
$ws.Range("K102").Value = 2941.6
$ws.Range("L102").Value = 3176
$ws.Range("M102").Value = -1319.6
$ws.Range("N102").Value = -6420
$ws.Range("H113").Value = 9308.786
$ws.Range("I113").Value = 12041
$ws.Range("J113").Value = 2478.25
$ws.Range("K113").Value = 12041
$ws.Range("L113").Value = 2478.25
$ws.Range("M113").Value = -9871
$ws.Range("N113").Value = -6818.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2275.682
$ws.Range("I40").Value = 2189.1177
$ws.Range("J40").Value = 2570
$ws.Range("K40").Value = 2189.1177
$ws.Range("L40").Value = 2570
$ws.Range("M40").Value = -2053.1177
$ws.Range("N40").Value = -2842
$ws.Range("H46").Value = 3388.077
$ws.Range("I46").Value = 423.16666
$ws.Range("J46").Value = 5929.4287
$ws.Range("K46").Value = 423.16666
$ws.Range("L46").Value = 5929.4287
$ws.Range("M46").Value = -235.16666
$ws.Range("N46").Value = -6305.4287
$ws.Range("H61").Value = 2375.5
$ws.Range("I61").Value = 2475.25
$ws.Range("J61").Value = 2176
$ws.Range("K61").Value = 2475.25
$ws.Range("L61").Value = 2176
$ws.Range("M61").Value = -2273.25
$ws.Range("N61").Value = -2580
$ws.Range("H113").Value = 2375.5
$ws.Range("I113").Value = 2475.25
$ws.Range("J113").Value = 2176
$ws.Range("K113").Value = 2475.25
$ws.Range("L113").Value = 2176
$ws.Range("M113").Value = -305.25
$ws.Range("N113").Value = -6516
$ws.Range("H132").Value = 4662.52
$ws.Range("I132").Value = 5864.9
$ws.Range("J132").Value = 3860.9333
$ws.Range("K132").Value = 17594.7
$ws.Range("L132").Value = 11582.7999
$ws.Range("M132").Value = -15064.7
$ws.Range("N132").Value = -16642.7999
$ws.Range("H136").Value = 3187.7058
$ws.Range("I136").Value = 2463.818
$ws.Range("J136").Value = 4514.8335
$ws.Range("K136").Value = 7391.454000000001
$ws.Range("L136").Value = 13544.5005
$ws.Range("M136").Value = -4841.454000000001
$ws.Range("N136").Value = -18644.5005
$ws.Range("H141").Value = 34761.383
$ws.Range("J141").Value = 34761.383
$ws.Range("L141").Value = 34761.383
$ws.Range("N141").Value = -45121.383

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1773.4445
$ws.Range("I81").Value = 1677.5
$ws.Range("J81").Value = 1850.2
$ws.Range("K81").Value = 3355
$ws.Range("L81").Value = 3700.4
$ws.Range("M81").Value = -2294
$ws.Range("N81").Value = -5822.4
$ws.Range("H84").Value = 1773.4445
$ws.Range("I84").Value = 1677.5
$ws.Range("J84").Value = 1850.2
$ws.Range("K84").Value = 16775
$ws.Range("L84").Value = 18502
$ws.Range("M84").Value = -11471
$ws.Range("N84").Value = -29110
$ws.Range("H132").Value = 2014.4
$ws.Range("I132").Value = 1268.375
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 3805.125
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -1275.125
$ws.Range("N132").Value = -20055.5
$ws.Range("H136").Value = 21574.215
$ws.Range("I136").Value = 55957.8
$ws.Range("J136").Value = 2472.2222
$ws.Range("K136").Value = 167873.4
$ws.Range("L136").Value = 7416.6666
$ws.Range("M136").Value = -165323.4
$ws.Range("N136").Value = -12516.6666
